$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("latest")

$ws.Range("B2").Value = 0.186985831112528
$ws.Range("C2").Value = 0.722203634391301
$ws.Range("D2").Value = 0.9341118818059831
$ws.Range("E2").Value = 0.9664946362013517
$ws.Range("F2").Value = 0.9757249869846261
$ws.Range("G2").Value = 18
$ws.Range("B3").Value = 0.2330935088154202
$ws.Range("C3").Value = 0.4998167590291169
$ws.Range("D3").Value = 0.3374744503351871
$ws.Range("E3").Value = 0.580925511864634
$ws.Range("F3").Value = 0.5484872224022262
$ws.Range("G3").Value = 17
$ws.Range("B4").Value = 0.2474022870896494
$ws.Range("C4").Value = 0.4562231751086744
$ws.Range("D4").Value = 0.282823943946844
$ws.Range("E4").Value = 0.531811944155868
$ws.Range("F4").Value = 0.4862000162233972
$ws.Range("G4").Value = 16
$ws.Range("B5").Value = 0.3331709995737225
$ws.Range("C5").Value = 0.5061089045982429
$ws.Range("D5").Value = 0.3514486812937936
$ws.Range("E5").Value = 0.5928310731513604
$ws.Range("F5").Value = 0.5075632610151456
$ws.Range("G5").Value = 15
$ws.Range("B6").Value = 0.3697548977873877
$ws.Range("C6").Value = 0.5076571467308241
$ws.Range("D6").Value = 0.3515889833503586
$ws.Range("E6").Value = 0.5929493935829251
$ws.Range("F6").Value = 0.4810392743262605
$ws.Range("G6").Value = 14
$ws.Range("B7").Value = 0.3550099532386083
$ws.Range("C7").Value = 0.5235044444458751
$ws.Range("D7").Value = 0.3739286097767974
$ws.Range("E7").Value = 0.611497023522435
$ws.Range("F7").Value = 0.518222527605833
$ws.Range("G7").Value = 13
$ws.Range("B8").Value = 0.4180599993501528
$ws.Range("C8").Value = 0.5377036271313899
$ws.Range("D8").Value = 0.3917271539151825
$ws.Range("E8").Value = 0.6258811020594747
$ws.Range("F8").Value = 0.4864935662755369
$ws.Range("G8").Value = 12
$ws.Range("B9").Value = 0.4373552979544756
$ws.Range("C9").Value = 0.5644811036562025
$ws.Range("D9").Value = 0.414358809081934
$ws.Range("E9").Value = 0.643707083293274
$ws.Range("F9").Value = 0.4953655899196011
$ws.Range("G9").Value = 11
$ws.Range("B10").Value = 0.4113383373941182
$ws.Range("C10").Value = 0.553516800113468
$ws.Range("D10").Value = 0.4096270470205588
$ws.Range("E10").Value = 0.640021130135997
$ws.Range("F10").Value = 0.5168578347523529
$ws.Range("G10").Value = 10
$ws.Range("B11").Value = 0.3463436447865475
$ws.Range("C11").Value = 0.5120086005850197
$ws.Range("D11").Value = 0.3587015581898463
$ws.Range("E11").Value = 0.5989169877285552
$ws.Range("F11").Value = 0.5182577473073058
$ws.Range("G11").Value = 9